# Angptl3-Itgav.xlsx was regenerated from the NATMI pipeline with updated
# per-cluster TPM expression values. This refreshes the ligand (G,H) and
# receptor (M,N) expression figures and recomputes every value that is
# derived from them: the ligand/receptor specificity columns (I,J,O,P)
# and the edge-weight columns (Q,R,S,T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies rows 2-21: 4 "Sending cluster" blocks of 5 rows each,
# every block cycling through the same 5 "Target cluster" values
# (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac) in that order.
$firstDataRow = 2
$lastDataRow = 21
$blockSize = 5

# --- New ligand average/total expression (TPM) values, one pair per
#     "Sending cluster" block (column G = average, column H = total) ---
$ligandG = @(1.363865666666667, 1.444951, 0.680678, 1.220628333333333)
$ligandH = @(4.091597, 4.334853, 2.042034, 3.661885)

# --- New receptor average/total expression (TPM) values, one pair per
#     "Target cluster" position within a block (column M = average,
#     column N = total) ---
$receptorM = @(22.495411, 82.64333833333332, 79.32606499999999, 14.467164, 72.34725666666667)
$receptorN = @(67.486233, 247.930015, 237.978195, 43.401492, 217.04177)

# 1) Write the refreshed raw expression values into G, H, M, N.
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $blockIndex = [math]::Floor(($row - $firstDataRow) / $blockSize)
    $posIndex = ($row - $firstDataRow) % $blockSize

    $ws.Range("G$row").Value2 = $ligandG[$blockIndex]
    $ws.Range("H$row").Value2 = $ligandH[$blockIndex]
    $ws.Range("M$row").Value2 = $receptorM[$posIndex]
    $ws.Range("N$row").Value2 = $receptorN[$posIndex]
}

# 2) Derived-specificity normalisers: sum of the average/total ligand
#    expression across the 4 sending clusters, and sum of the
#    average/total receptor expression across the 5 target clusters.
$sumLigandG = 0
$sumLigandH = 0
foreach ($v in $ligandG) { $sumLigandG += $v }
foreach ($v in $ligandH) { $sumLigandH += $v }

$sumReceptorM = 0
$sumReceptorN = 0
foreach ($v in $receptorM) { $sumReceptorM += $v }
foreach ($v in $receptorN) { $sumReceptorN += $v }

# 3) Recompute ligand specificity (I,J), receptor specificity (O,P),
#    and the edge-weight columns (Q,R,S,T) for every row.
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $blockIndex = [math]::Floor(($row - $firstDataRow) / $blockSize)
    $posIndex = ($row - $firstDataRow) % $blockSize

    $g = $ligandG[$blockIndex]
    $h = $ligandH[$blockIndex]
    $m = $receptorM[$posIndex]
    $n = $receptorN[$posIndex]

    $i = $g / $sumLigandG
    $j = $h / $sumLigandH
    $o = $m / $sumReceptorM
    $p = $n / $sumReceptorN

    $ws.Range("I$row").Value2 = $i
    $ws.Range("J$row").Value2 = $j
    $ws.Range("O$row").Value2 = $o
    $ws.Range("P$row").Value2 = $p

    $ws.Range("Q$row").Value2 = $g * $m
    $ws.Range("R$row").Value2 = $h * $n
    $ws.Range("S$row").Value2 = $i * $o
    $ws.Range("T$row").Value2 = $j * $p
}
